$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "Tumor" query text (row 3, column B on the startup sheet) is corrected:
# the RETURN clause now projects the sample's own tumor-status property
# instead of the aggregated WITH-alias, and the final ORDER BY line is
# re-indented by two spaces.
$tumorQuery = "MATCH (s:study)<--(p:participant)<--(samp:sample)`r`n" +
  "WHERE s.study_name in [`"Human Tumor Atlas Network (HTAN) primary sequencing data`"]`r`n" +
  "WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`r`n" +
  "RETURN  `r`n" +
  " coalesce(samp.sample_id, '') as ``Sample ID``,`r`n" +
  " coalesce(p.participant_id,'') as ``Participant ID``,`r`n" +
  " coalesce(s.study_name, '') as ``Study Name``,`r`n" +
  " coalesce(s.phs_accession,'') as ``Accession``,`r`n" +
  " coalesce(samp.sample_tumor_status,'') as ``Tumor``,`r`n" +
  "coalesce(samp.sample_type,'') as ``Analyte Type```r`n" +
  "  ORDER By samp.sample_id LIMIT 100"

$ws.Range("B3").Value = $tumorQuery
